$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 345, shifting all existing rows (345..687) down by one.
$ws.Rows(345).Insert()

# Populate the newly inserted row with the new entity/category/count values.
$ws.Range("A345").Value = "Section Ecosystem Functions"
$ws.Range("B345").Value = "ENVIRONMENT"
$ws.Range("C345").Value = 1
